$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '41.231.41'
$ws.Range('E2').Value = '  -1.09%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.430.58'
$ws.Range('E3').Value = '  -1.91%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  +0.20%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '316.80'
$ws.Range('E5').Value = '  -0.70%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '88.70'
$ws.Range('E6').Value = '  -4.71%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.541'
$ws.Range('E7').Value = '  -2.54%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.495'
$ws.Range('E9').Value = '  -4.71%  '
$ws.Range('B10').Value = 'Dogecoin'
$ws.Range('C10').Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0833'
$ws.Range('E10').Value = '  -6.22%  '
$ws.Range('B11').Value = 'Avalanche'
$ws.Range('C11').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '32.04'
$ws.Range('E11').Value = '  -3.14%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.108'
$ws.Range('E12').Value = '  -2.85%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '2.808.47'
$ws.Range('E13').Value = '  -1.75%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.70'
$ws.Range('E14').Value = '  -3.65%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '15.53'
$ws.Range('E15').Value = '  -1.21%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '2.425.02'
$ws.Range('E16').Value = '  -2.49%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.771'
$ws.Range('E17').Value = '  -2.76%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '41.145.83'
$ws.Range('E18').Value = '  -1.18%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.0₃0923'
$ws.Range('E19').Value = '  -4.17%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '6.22'
$ws.Range('E20').Value = '  -4.38%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '71.58'
$ws.Range('E21').Value = '  +0.02%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '11.01'
$ws.Range('E22').Value = '  -4.69%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '235.28'
$ws.Range('E23').Value = '  -2.70%  '
$ws.Range('E24').Value = '  -2.62%  '
$ws.Range('E25').Value = '  -0.02%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.87'
$ws.Range('E26').Value = '  -3.02%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '23.93'
$ws.Range('E27').Value = '  -3.97%  '
$ws.Range('E28').Value = '  -3.40%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '9.53'
$ws.Range('E29').Value = '  -3.72%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '34.52'
$ws.Range('E30').Value = '  -5.79%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '157.13'
$ws.Range('E31').Value = '  +0.18%  '
$ws.Range('B32').Value = 'FirstDigitalUSD'
$ws.Range('C32').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.00'
$ws.Range('E32').Value = '  +0.11%  '
$ws.Range('B33').Value = 'Filecoin'
$ws.Range('C33').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.25'
$ws.Range('E33').Value = '  -5.24%  '
$ws.Range('E34').Value = '  -1.79%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.0742'
$ws.Range('E35').Value = '  -3.93%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.89'
$ws.Range('E36').Value = '  -1.42%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '16.51'
$ws.Range('E37').Value = '  -5.96%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.114'
$ws.Range('E38').Value = '  -0.97%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.76'
$ws.Range('E39').Value = '  -4.00%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.0995'
$ws.Range('E40').Value = '  -4.01%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '3.86'
$ws.Range('E41').Value = '  -4.05%  '
$ws.Range('E42').Value = '  -7.11%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '1.979.82'
$ws.Range('E43').Value = '  -0.24%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0274'
$ws.Range('E44').Value = '  -4.25%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '18.15'
$ws.Range('E45').Value = '  -7.44%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.86'
$ws.Range('E46').Value = '  -5.81%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '9.51'
$ws.Range('E47').Value = '  +3.20%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.672.72'
$ws.Range('E48').Value = '  -1.49%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '95.15'
$ws.Range('E49').Value = '  -2.68%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '73.07'
$ws.Range('E50').Value = '  -1.56%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '51.91'
$ws.Range('E51').Value = '  -1.80%  '
